$d = $word.ActiveDocument

# Locate the "footer" block that should be removed: the blank paragraph
# right after "LOB1053: Física III (Requisito fraco)", the
# "Ver no Jupiter ..." paragraph, and the "© 2020 ..." paragraph.
# They are found by searching for the distinctive copyright text and the
# "Ver no Jupiter" text, then widened to cover the blank paragraph that
# immediately precedes them.

$find = $d.Content.Duplicate
$find.Find.ClearFormatting()
$ok = $find.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not find the 'Ver no Jupiter' paragraph"
}

$jupiterPara = $find.Paragraphs(1)
$blankPara = $jupiterPara.Previous(1)
$copyrightPara = $jupiterPara.Next(1)

# Sanity-check we grabbed the right paragraphs before deleting anything.
if ($blankPara.Range.Text.Trim() -ne "") {
    throw "Expected a blank paragraph before 'Ver no Jupiter', found: $($blankPara.Range.Text)"
}
if ($copyrightPara.Range.Text -notmatch "Contact: luizeleno@usp.br") {
    throw "Expected the copyright paragraph after 'Ver no Jupiter', found: $($copyrightPara.Range.Text)"
}

$start = $blankPara.Range.Start
$end = $copyrightPara.Range.End

$killRange = $d.Range($start, $end)
$killRange.Delete()
